$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.571.44"
$ws.Range("E2").Value = "  +0.18%  "

$ws.Range("D3").Value = "'1.922.04"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "'326.14"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  +0.35%  "

$ws.Range("D7").Value = "'0.4818"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").Value = "'0.4066"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").Value = "'0.08226"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("E10").Value = "  -0.24%  "

$ws.Range("D11").Value = "'23.61"
$ws.Range("E11").Value = "  +1.08%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "'1.911.31"
$ws.Range("E12").Value = "  -2.51%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.083"
$ws.Range("E13").Value = "  +1.35%  "

$ws.Range("D15").Value = "'91.63"
$ws.Range("E15").Value = "  +1.47%  "

$ws.Range("D16").Value = "'0.06868"
$ws.Range("E16").Value = "  +1.19%  "

$ws.Range("E17").Value = "  +0.27%  "

$ws.Range("E18").Value = "  -0.04%  "

$ws.Range("E19").Value = "  -0.27%  "

$ws.Range("E20").Value = "  +0.28%  "

$ws.Range("D21").Value = "'29.584.05"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("D22").Value = "'5.684"
$ws.Range("E22").Value = "  +1.16%  "

$ws.Range("D24").Value = "'2.184"
$ws.Range("E24").Value = "  -0.03%  "

$ws.Range("D25").Value = "'2.137.16"
$ws.Range("E25").Value = "  -2.12%  "

$ws.Range("D26").Value = "'155.86"
$ws.Range("E26").Value = "  +0.36%  "

$ws.Range("D27").Value = "'6.457"
$ws.Range("E27").Value = "  -0.18%  "

$ws.Range("E28").Value = "  -0.35%  "

$ws.Range("D29").Value = "'2.095"
$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "'120.58"
$ws.Range("E30").Value = "  +0.85%  "

$ws.Range("D31").Value = "'1.015"
$ws.Range("E31").Value = "  -1.49%  "

$ws.Range("D32").Value = "'0.09635"
$ws.Range("E32").Value = "  +0.76%  "

$ws.Range("D33").Value = "'5.620"
$ws.Range("E33").Value = "  +1.82%  "

$ws.Range("D34").Value = "'3.554"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").Value = "'1.379"
$ws.Range("E35").Value = "  -1.07%  "

$ws.Range("D36").Value = "'0.06362"
$ws.Range("E36").Value = "  +4.33%  "

$ws.Range("E37").Value = "  +1.21%  "

$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("D39").Value = "'0.5961"
$ws.Range("E39").Value = "  +0.50%  "

$ws.Range("D40").Value = "'10.77"
$ws.Range("E40").Value = "  -0.05%  "

$ws.Range("E41").Value = "  +0.29%  "

$ws.Range("D42").Value = "'7.885"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").Value = "'2.455"
$ws.Range("E44").Value = "  -1.08%  "

$ws.Range("D45").Value = "'1.271"
$ws.Range("E45").Value = "  -0.99%  "

$ws.Range("D46").Value = "'12.44"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").Value = "'0.07485"
$ws.Range("E47").Value = "  -3.03%  "

$ws.Range("D48").Value = "'0.5566"
$ws.Range("E48").Value = "  -0.10%  "

$ws.Range("E49").Value = "  +0.36%  "

$ws.Range("D50").Value = "'119.09"
$ws.Range("E50").Value = "  +2.95%  "

$ws.Range("E51").Value = "  +3.34%  "
